# Add a new daily snapshot row (id=7) to the "Daily APR" data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 8
$ws.Cells.Item($row, 1).Value = 7
$ws.Cells.Item($row, 2).Value = "2025-09-11T09:31"
$ws.Cells.Item($row, 3).Value = 1.6792199976262983
